# EPBDS-13025 The empty array elements are trimmed
#
# The "Value2" column (column D, rows 38-58) on Sheet1 holds the expected
# textual results of sorting arrays that contain a null/empty element.
# Previously the trailing empty element produced a stray trailing comma
# (e.g. "1, 2, ,"). After the fix the empty element is trimmed from the
# joined string, so the trailing comma disappears (e.g. "1, 2, ").
#
# Trim the trailing "," from every non-empty cell in D38:D58.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($row = 38; $row -le 58; $row++) {
    $cell = $ws.Cells.Item($row, 4)
    $current = $cell.Value2
    if ($current -ne $null -and $current.Length -gt 0 -and $current.EndsWith(",")) {
        $trimmed = $current.Substring(0, $current.Length - 1)
        # Leading apostrophe preserves the existing quote-prefixed text style
        # (the cell already stores text, not a number, so re-apply the same
        # "entered as text" marker Excel uses).
        $cell.Value = "'" + $trimmed
    }
}
